# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# per latest scrape, matching the commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "36.767.17"
$ws.Range("E2").Value = "  +4.27%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.924.14"
$ws.Range("E3").Value = "  +2.47%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5: BNB
$ws.Range("D5").Value = "'249.71"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.700"
$ws.Range("E6").Value = "  +2.91%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8: Solana
$ws.Range("D8").Value = "'44.53"
$ws.Range("E8").Value = "  +2.34%  "

# Row 9: OKB
$ws.Range("D9").Value = "'58.62"
$ws.Range("E9").Value = "  +9.46%  "

# Row 10: Cardano
$ws.Range("E10").Value = "  +4.21%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.0765"
$ws.Range("E11").Value = "  +3.92%  "

# Row 12: TRON
$ws.Range("D12").Value = "'0.0999"
$ws.Range("E12").Value = "  +2.75%  "

# Row 13: Chainlink
$ws.Range("D13").Value = "'14.66"
$ws.Range("E13").Value = "  +8.88%  "

# Row 14: Polygon
$ws.Range("D14").Value = "'0.804"
$ws.Range("E14").Value = "  +5.33%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.201.34"
$ws.Range("E15").Value = "  +2.45%  "

# Row 16: Polkadot
$ws.Range("E16").Value = "  +4.87%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "1.922.67"
$ws.Range("E17").Value = "  +2.44%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "36.693.05"
$ws.Range("E18").Value = "  +4.00%  "

# Row 19: Litecoin
$ws.Range("D19").Value = "'74.35"
$ws.Range("E19").Value = "  +2.22%  "

# Row 20: ShibaInu
$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("E20").Value = "  +5.14%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'252.34"
$ws.Range("E21").Value = "  +3.62%  "

# Row 22: Avalanche
$ws.Range("D22").Value = "'13.32"
$ws.Range("E22").Value = "  +4.42%  "

# Row 23: Uniswap
$ws.Range("D23").Value = "'5.22"
$ws.Range("E23").Value = "  +5.64%  "

# Row 24: Toncoin
$ws.Range("D24").Value = "'2.68"
$ws.Range("E24").Value = "  +2.32%  "

# Row 25: Dai
$ws.Range("E25").Value = "  +0.02%  "

# Row 26: PancakeSwap
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  +0.72%  "

# Row 27: Monero
$ws.Range("D27").Value = "'168.49"
$ws.Range("E27").Value = "  +1.72%  "

# Row 28: Cosmos
$ws.Range("D28").Value = "'8.86"
$ws.Range("E28").Value = "  +4.43%  "

# Row 29: EthereumClassic
$ws.Range("E29").Value = "  +3.43%  "

# Row 30: Stellar
$ws.Range("D30").Value = "'0.130"
$ws.Range("E30").Value = "  +2.59%  "

# Row 31: Filecoin
$ws.Range("E31").Value = "  +6.49%  "

# Row 32: Hedera
$ws.Range("D32").Value = "'0.0621"
$ws.Range("E32").Value = "  +5.26%  "

# Row 33: WEMIXToken
$ws.Range("D33").Value = "'1.98"
$ws.Range("E33").Value = "  -1.68%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +5.82%  "

# Row 35: BinanceUSD
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36: Kaspa
$ws.Range("D36").Value = "'0.0865"
$ws.Range("E36").Value = "  +19.28%  "

# Row 37: TrustWalletToken
$ws.Range("D37").Value = "'1.52"
$ws.Range("E37").Value = "  -10.21%  "

# Row 38: ImmutableX
$ws.Range("D38").Value = "'0.901"
$ws.Range("E38").Value = "  +7.67%  "

# Row 39: Gas
$ws.Range("D39").Value = "'17.89"
$ws.Range("E39").Value = "  +50.04%  "

# Row 40: LidoDAOToken
$ws.Range("E40").Value = "  +4.53%  "

# Row 41: Aave
$ws.Range("D41").Value = "'106.63"
$ws.Range("E41").Value = "  +11.35%  "

# Row 42: VeChain
$ws.Range("E42").Value = "  +5.51%  "

# Row 43: InjectiveProtocol
$ws.Range("E43").Value = "  -2.02%  "

# Row 44: ARBITRUM
$ws.Range("E44").Value = "  +3.91%  "

# Row 45: Maker
$ws.Range("D45").Value = "1.339.36"
$ws.Range("E45").Value = "  +2.92%  "

# Row 46: HuobiToken
$ws.Range("D46").Value = "'2.58"
$ws.Range("E46").Value = "  +8.18%  "

# Row 47: RenderToken
$ws.Range("E47").Value = "  +1.43%  "

# Row 48: Cronos
$ws.Range("D48").Value = "'0.0816"
$ws.Range("E48").Value = "  +2.42%  "

# Row 49: MXToken
$ws.Range("E49").Value = "  +2.91%  "

# Row 50: FraxShare
$ws.Range("D50").Value = "'6.46"
$ws.Range("E50").Value = "  +4.19%  "

# Row 51: MultiversX
$ws.Range("D51").Value = "'43.36"
$ws.Range("E51").Value = "  +3.54%  "
